$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the data (row 2), pushing existing data down.
$ws.Rows.Item(2).Insert()

# The freshly inserted row doesn't inherit the plain "data row" formatting
# used throughout the table (centered, unbordered, non-bold). Match it.
$rng = $ws.Range("A2:F2")
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$rng.Borders.LineStyle = -4142
$ws.Range("D2").NumberFormat = "0.000"

# New day's record: 05-11-2025, still referencing the 01-11-2025 circular.
# Leading apostrophes keep the dd-mm-yyyy strings as literal text instead of
# being auto-converted to date serials.
$ws.Range("A2").Value = "'05-11-2025"
$ws.Range("B2").Value = "ALUMINIUM INGOT"
$ws.Range("C2").Value = "IE07"
$ws.Range("D2").Value = 297.15
$ws.Range("E2").Value = "'01-11-2025"
$url = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf"
$ws.Range("F2").Value = $url
$ws.Hyperlinks.Add($ws.Range("F2"), $url)

# Adding the hyperlink auto-applies Excel's blue/underlined "Hyperlink" style;
# reset F2 back to the plain look the rest of the column uses.
$ws.Range("F2").Font.Bold = $false
$ws.Range("F2").Font.Underline = $false
$ws.Range("F2").Font.ThemeColor = 1
$ws.Range("F2").Font.ColorIndex = -4105
$ws.Range("F2").HorizontalAlignment = -4108
$ws.Range("F2").VerticalAlignment = -4108
